$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: étude number
$ws.Range("H2").Value = "N° 24001"

# Date
$ws.Range("I3").Value = "21 July 2024"

# Interlocuteur name casing
$ws.Range("G4").Value = "Antony Feord"

# Address lines
$ws.Range("G6").Value = "4 rue des M"
$ws.Range("G8").Value = "77420 France"

# N° SS now a text "0" instead of a long numeric SS number
$ws.Range("H10").Value = "'0"

# Etude reference (must stay text, not be parsed as 24e01 = 240 numeric)
$ws.Range("C13").Value = "'24e01"

# Rétribution brute
$ws.Range("I13").Value = 560

# Nb de Jours-Etude Homme
$ws.Range("I14").Value = 2

# Taux AT/MP
$ws.Range("F23").Value = 0.66
